$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("M2").Value = 2.843949
$ws.Range("N2").Value = 8.531846999999999
$ws.Range("O2").Value = 0.4976240243095911
$ws.Range("P2").Value = 0.4976240243095912
$ws.Range("Q2").Value = 34.28859440511599
$ws.Range("R2").Value = 308.5973496460439
$ws.Range("S2").Value = 0.03138923996367713
$ws.Range("T2").Value = 0.03138923996367714

$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("O3").Value = 0.4403664892852895
$ws.Range("P3").Value = 0.4403664892852897
$ws.Range("Q3").Value = 30.343285699796
$ws.Range("R3").Value = 273.089571298164
$ws.Range("S3").Value = 0.02777753631030147
$ws.Range("T3").Value = 0.02777753631030148

$ws.Range("G4").Value = 12.056684
$ws.Range("H4").Value = 36.170052
$ws.Range("I4").Value = 0.06307822458376462
$ws.Range("J4").Value = 0.06307822458376462
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("O4").Value = 0.06200948640511928
$ws.Range("P4").Value = 0.0620094864051193
$ws.Range("Q4").Value = 4.272740110497332
$ws.Range("R4").Value = 38.45466099447599
$ws.Range("S4").Value = 0.003911448309786013
$ws.Range("T4").Value = 0.003911448309786014

$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("M5").Value = 2.843949
$ws.Range("N5").Value = 8.531846999999999
$ws.Range("O5").Value = 0.4976240243095911
$ws.Range("P5").Value = 0.4976240243095912
$ws.Range("Q5").Value = 71.498680923666
$ws.Range("R5").Value = 643.4881283129939
$ws.Range("S5").Value = 0.0654529382593904
$ws.Range("T5").Value = 0.06545293825939041

$ws.Range("I6").Value = 0.1315309049843414
$ws.Range("J6").Value = 0.1315309049843414
$ws.Range("O6").Value = 0.4403664892852895
$ws.Range("P6").Value = 0.4403664892852897
$ws.Range("S6").Value = 0.0579218028604714
$ws.Range("T6").Value = 0.05792180286047141

$ws.Range("I7").Value = 0.1315309049843414
$ws.Range("J7").Value = 0.1315309049843414
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("O7").Value = 0.06200948640511928
$ws.Range("P7").Value = 0.0620094864051193
$ws.Range("Q7").Value = 8.909530621780666
$ws.Range("R7").Value = 80.18577559602599
$ws.Range("S7").Value = 0.008156163864479551
$ws.Range("T7").Value = 0.008156163864479555

$ws.Range("H8").Value = 461.8238680000001
$ws.Range("I8").Value = 0.8053908704318941
$ws.Range("J8").Value = 0.8053908704318941
$ws.Range("M8").Value = 2.843949
$ws.Range("N8").Value = 8.531846999999999
$ws.Range("O8").Value = 0.4976240243095911
$ws.Range("P8").Value = 0.4976240243095912
$ws.Range("Q8").Value = 437.801175858244
$ws.Range("R8").Value = 3940.210582724196
$ws.Range("S8").Value = 0.4007818460865236
$ws.Range("T8").Value = 0.4007818460865237

$ws.Range("H9").Value = 461.8238680000001
$ws.Range("I9").Value = 0.8053908704318941
$ws.Range("J9").Value = 0.8053908704318941
$ws.Range("O9").Value = 0.4403664892852895
$ws.Range("P9").Value = 0.4403664892852897
$ws.Range("R9").Value = 3486.842709747277
$ws.Range("S9").Value = 0.3546671501145167
$ws.Range("T9").Value = 0.3546671501145168

$ws.Range("H10").Value = 461.8238680000001
$ws.Range("I10").Value = 0.8053908704318941
$ws.Range("J10").Value = 0.8053908704318941
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("O10").Value = 0.06200948640511928
$ws.Range("P10").Value = 0.0620094864051193
$ws.Range("Q10").Value = 54.55489433049821
$ws.Range("S10").Value = 0.04994187423085372
$ws.Range("T10").Value = 0.04994187423085374
